$wb = $excel.ActiveWorkbook

# Update "想去人数" (number of people wanting to go) values on the
# "展览" and "全部类型" worksheets.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 5888
    $ws.Range("F5").Value = 985
    $ws.Range("F6").Value = 74
}
